$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two store names between row 4 and row 5
$ws.Range("A4").Value = "Bibi Cell Manauara"
$ws.Range("A5").Value = "Bibi Cell Ponta Negra"

# Swap the daily sales figures (columns B:H) between row 4 and row 5,
# and add the new day's figures in column I for both rows.
$ws.Range("B4").Value = 3340
$ws.Range("C4").Value = 1519
$ws.Range("D4").Value = 2934
$ws.Range("E4").Value = 1819
$ws.Range("F4").Value = 2503
$ws.Range("G4").Value = 2892
$ws.Range("H4").Value = 4208.4
$ws.Range("I4").Value = 3329.9

$ws.Range("B5").Value = 1800.01
$ws.Range("C5").Value = 4670
$ws.Range("D5").Value = 1748.51
$ws.Range("E5").Value = 5592
$ws.Range("F5").Value = 3002
$ws.Range("G5").Value = 823
$ws.Range("H5").Value = 3138.5
$ws.Range("I5").Value = 1613

# Update row totals (column AG) for rows 4 and 5
$ws.Range("AG4").Value = 22545.3
$ws.Range("AG5").Value = 22387.02

# Update the combined total row (row 6): new day column I and grand total AG
$ws.Range("I6").Value = 4942.9
$ws.Range("AG6").Value = 155464.58
